$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I1 = "I0" and J1 = "IF", formatted like the other header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: add numeric values for columns I and J
$data = @{
    2 = @(6, 6)
    3 = @(7, 7)
    4 = @(8, 8)
    5 = @(6, 6)
    6 = @(6, 6)
    7 = @(7, 7)
    8 = @(2, 3)
    9 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
